$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Scanner" to "Session"
$ws.Name = "Session"

# Remove the logged scan entry in row 2 (student 755263 / backup@backdoor.com)
# so only the header row remains.
$ws.Rows(2).Delete()
